# Updated cryptos list on Mon Jan  1 11:48:11 UTC 2024 with GitHub Actions
#
# The "Price" column (D) holds numeric-looking text (e.g. "309.56",
# "42.981.47") that must stay stored as text, exactly like the source
# workbook. A leading apostrophe forces Excel to keep the literal text
# instead of silently coercing it to a number (which would also strip
# meaningful trailing zeros, e.g. "36.20" -> 36.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indexes: A=1, B=2, C=3, D=4, E=5

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "'42.981.47"
$ws.Cells.Item(2, 5).Value = "  +0.54%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "'2.300.40"
$ws.Cells.Item(3, 5).Value = "  -0.51%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.14%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "'309.56"
$ws.Cells.Item(5, 5).Value = "  -2.76%  "

# Row 6 - Solana
$ws.Cells.Item(6, 4).Value = "'104.63"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "

# Row 7 - XRP
$ws.Cells.Item(7, 4).Value = "'0.618"
$ws.Cells.Item(7, 5).Value = "  -1.91%  "

# Row 8 - USDC
$ws.Cells.Item(8, 5).Value = "  -0.07%  "

# Row 9 - Cardano
$ws.Cells.Item(9, 5).Value = "  -0.67%  "

# Row 10 - Avalanche
$ws.Cells.Item(10, 4).Value = "'39.46"
$ws.Cells.Item(10, 5).Value = "  -1.47%  "

# Row 11 - Dogecoin
$ws.Cells.Item(11, 4).Value = "'0.0902"
$ws.Cells.Item(11, 5).Value = "  -0.50%  "

# Row 12 - Polkadot
$ws.Cells.Item(12, 4).Value = "'8.25"
$ws.Cells.Item(12, 5).Value = "  -3.31%  "

# Row 13 - TRON
$ws.Cells.Item(13, 5).Value = "  +0.57%  "

# Row 14 - Polygon
$ws.Cells.Item(14, 4).Value = "'0.989"
$ws.Cells.Item(14, 5).Value = "  +1.15%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = "'2.781.75"
$ws.Cells.Item(15, 5).Value = "  +4.43%  "

# Row 16 - Chainlink
$ws.Cells.Item(16, 4).Value = "'15.34"
$ws.Cells.Item(16, 5).Value = "  -0.53%  "

# Row 17 - WrappedEther
$ws.Cells.Item(17, 4).Value = "'2.293.18"
$ws.Cells.Item(17, 5).Value = "  -1.39%  "

# Row 18 - WrappedBTC
$ws.Cells.Item(18, 4).Value = "'42.760.41"
$ws.Cells.Item(18, 5).Value = "  +0.14%  "

# Row 19 - Uniswap
$ws.Cells.Item(19, 4).Value = "'7.31"
$ws.Cells.Item(19, 5).Value = "  -4.24%  "

# Row 20 - ShibaInu
$ws.Cells.Item(20, 5).Value = "  -1.26%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Cells.Item(21, 4).Value = "'13.54"
$ws.Cells.Item(21, 5).Value = "  +0.52%  "

# Row 22 - Litecoin
$ws.Cells.Item(22, 4).Value = "'73.25"
$ws.Cells.Item(22, 5).Value = "  -1.04%  "

# Row 23 - PancakeSwap
$ws.Cells.Item(23, 4).Value = "'3.43"
$ws.Cells.Item(23, 5).Value = "  -3.73%  "

# Row 24 - BitcoinCash
$ws.Cells.Item(24, 4).Value = "'267.08"
$ws.Cells.Item(24, 5).Value = "  -0.75%  "

# Row 25 - ImmutableX
$ws.Cells.Item(25, 5).Value = "  -1.80%  "

# Row 26 - Dai
$ws.Cells.Item(26, 5).Value = "  +0.49%  "

# Row 27 - Filecoin
$ws.Cells.Item(27, 4).Value = "'7.38"
$ws.Cells.Item(27, 5).Value = "  +18.07%  "

# Row 28 - Cosmos
$ws.Cells.Item(28, 4).Value = "'10.89"
$ws.Cells.Item(28, 5).Value = "  +0.03%  "

# Row 29 - Toncoin
$ws.Cells.Item(29, 5).Value = "  -1.21%  "

# Row 30 - EthereumClassic
$ws.Cells.Item(30, 4).Value = "'22.24"
$ws.Cells.Item(30, 5).Value = "  -2.07%  "

# Row 31 - InjectiveProtocol
$ws.Cells.Item(31, 4).Value = "'36.20"
$ws.Cells.Item(31, 5).Value = "  -4.48%  "

# Row 32 - Monero
$ws.Cells.Item(32, 4).Value = "'165.34"
$ws.Cells.Item(32, 5).Value = "  -0.17%  "

# Row 33 - Hedera
$ws.Cells.Item(33, 5).Value = "  -3.84%  "

# Row 34 - WEMIXToken
$ws.Cells.Item(34, 5).Value = "  +2.11%  "

# Row 35 - Stellar
$ws.Cells.Item(35, 4).Value = "'0.130"
$ws.Cells.Item(35, 5).Value = "  -1.42%  "

# Row 36 - Kaspa
$ws.Cells.Item(36, 5).Value = "  -3.66%  "

# Row 37 - RenderToken
$ws.Cells.Item(37, 4).Value = "'4.55"
$ws.Cells.Item(37, 5).Value = "  -1.53%  "

# Row 38 - VeChain
$ws.Cells.Item(38, 5).Value = "  -1.77%  "

# Row 39 - LidoDAOToken
$ws.Cells.Item(39, 5).Value = "  +1.66%  "

# Row 40 - NEARProtocol
$ws.Cells.Item(40, 5).Value = "  -2.52%  "

# Row 41 - BitcoinSV
$ws.Cells.Item(41, 4).Value = "'109.63"
$ws.Cells.Item(41, 5).Value = "  +12.05%  "

# Row 42 - ARBITRUM
$ws.Cells.Item(42, 5).Value = "  -3.94%  "

# Row 43 - MultiversX
$ws.Cells.Item(43, 4).Value = "'71.18"
$ws.Cells.Item(43, 5).Value = "  +1.15%  "

# Rows 44/45 - Algorand and FirstDigitalUSD swap ranking positions
$ws.Cells.Item(44, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(44, 4).Value = "'1.01"
$ws.Cells.Item(44, 5).Value = "  +0.28%  "

$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).Value = "'0.226"
$ws.Cells.Item(45, 5).Value = "  +0.00%  "

# Row 46 - Celestia
$ws.Cells.Item(46, 4).Value = "'12.23"
$ws.Cells.Item(46, 5).Value = "  -1.27%  "

# Row 47 - Maker
$ws.Cells.Item(47, 4).Value = "'1.732.92"
$ws.Cells.Item(47, 5).Value = "  +6.65%  "

# Row 48 - Aave
$ws.Cells.Item(48, 4).Value = "'110.59"
$ws.Cells.Item(48, 5).Value = "  -5.09%  "

# Row 49 - ordi
$ws.Cells.Item(49, 4).Value = "'77.04"
$ws.Cells.Item(49, 5).Value = "  -6.45%  "

# Row 50 - FraxShare
$ws.Cells.Item(50, 4).Value = "'8.64"
$ws.Cells.Item(50, 5).Value = "  -2.92%  "

# Row 51 - THORChain
$ws.Cells.Item(51, 4).Value = "'5.12"
$ws.Cells.Item(51, 5).Value = "  -3.34%  "
